$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column B (Number) to text format before writing numeric-looking strings,
# so Excel does not auto-convert them to numbers; formats are cleared afterward
# so the written cells keep the default (unstyled) appearance.
$ws.Range("B86:B101").NumberFormat = "@"

$ws.Cells.Item(86, 1).Value = "2026-02-16 11:44:54"
$ws.Cells.Item(86, 2).Value = "237671823369"
$ws.Cells.Item(86, 3).Value = "MFS ENTREE COLLEGE MALANGUE"
$ws.Cells.Item(86, 4).Value = 110633

$ws.Cells.Item(87, 1).Value = "2026-02-16 14:22:30"
$ws.Cells.Item(87, 2).Value = "237672128028"
$ws.Cells.Item(87, 3).Value = "CAROLINE WAKO DJAMNOU"
$ws.Cells.Item(87, 4).Value = 22483

$ws.Cells.Item(88, 1).Value = "2026-02-16 13:58:35"
$ws.Cells.Item(88, 2).Value = "237672277367"
$ws.Cells.Item(88, 3).Value = "TOP MOBIL KM5 LTDLA_POLAS_BTQ_KM5"
$ws.Cells.Item(88, 4).Value = 726946

$ws.Cells.Item(89, 1).Value = "2026-02-16 16:20:18"
$ws.Cells.Item(89, 2).Value = "237674853971"
$ws.Cells.Item(89, 3).Value = "NJOSSEU TCHOUNZOU TOP MOBILE"
$ws.Cells.Item(89, 4).Value = 129189

$ws.Cells.Item(90, 1).Value = "2026-02-16 17:06:21"
$ws.Cells.Item(90, 2).Value = "237674884705"
$ws.Cells.Item(90, 3).Value = "manuela verna yetna baaga"
$ws.Cells.Item(90, 4).Value = 15105

$ws.Cells.Item(91, 1).Value = "2026-02-16 15:21:25"
$ws.Cells.Item(91, 2).Value = "237675779272"
$ws.Cells.Item(91, 3).Value = "RODES NGWEM KEMAYOU"
$ws.Cells.Item(91, 4).Value = 28759

$ws.Cells.Item(92, 1).Value = "2026-02-16 14:23:53"
$ws.Cells.Item(92, 2).Value = "237677304210"
$ws.Cells.Item(92, 3).Value = "FERDINAND NKWELLE NGOME"
$ws.Cells.Item(92, 4).Value = 137078

$ws.Cells.Item(93, 1).Value = "2026-02-06 09:50:11"
$ws.Cells.Item(93, 2).Value = "237678267353"
$ws.Cells.Item(93, 3).Value = "LA NEGRESSE SARL EMBOLA BELTUS MBU"
$ws.Cells.Item(93, 4).Value = 0

$ws.Cells.Item(94, 1).Value = "2026-02-16 14:32:56"
$ws.Cells.Item(94, 2).Value = "237678370615"
$ws.Cells.Item(94, 3).Value = "ESSEN ONGOLONG BERTHE HORTENSE ETS MOBILE FINANCIAL SERVICES MFS"
$ws.Cells.Item(94, 4).Value = 239366

$ws.Cells.Item(95, 1).Value = "2026-02-16 13:21:38"
$ws.Cells.Item(95, 2).Value = "237678836319"
$ws.Cells.Item(95, 3).Value = "KAMDOM DOMINIQUE STEPHANIE ETS MOBILE FINANCIAL SERVICES MFS"
$ws.Cells.Item(95, 4).Value = 94879

$ws.Cells.Item(96, 1).Value = "2026-02-16 13:41:10"
$ws.Cells.Item(96, 2).Value = "237678922502"
$ws.Cells.Item(96, 3).Value = "NWOAGA TCHAMDJOU EPSE KAMSEU EMILINE ETS LE CONTENT"
$ws.Cells.Item(96, 4).Value = 1196062

$ws.Cells.Item(97, 1).Value = "2026-02-16 15:24:23"
$ws.Cells.Item(97, 2).Value = "237679884264"
$ws.Cells.Item(97, 3).Value = "MFS CICAM"
$ws.Cells.Item(97, 4).Value = 66267

$ws.Cells.Item(98, 1).Value = "2026-02-16 12:52:54"
$ws.Cells.Item(98, 2).Value = "237681019523"
$ws.Cells.Item(98, 3).Value = "ETS MOULAY RIPERT AND COMPANY"
$ws.Cells.Item(98, 4).Value = 132723

$ws.Cells.Item(99, 1).Value = "2026-02-16 16:55:47"
$ws.Cells.Item(99, 2).Value = "237681125655"
$ws.Cells.Item(99, 3).Value = "EMENGUE PICHOU ROMEO KAMILAH CONNECTION GROUP"
$ws.Cells.Item(99, 4).Value = 498608

$ws.Cells.Item(100, 1).Value = "2026-02-16 14:07:41"
$ws.Cells.Item(100, 2).Value = "237681240793"
$ws.Cells.Item(100, 3).Value = "MBANE EMILIE FRANCOISE ETS MOBILE FINANCIAL SERVICES MFS"
$ws.Cells.Item(100, 4).Value = 15143

$ws.Cells.Item(101, 1).Value = "2026-02-16 13:18:11"
$ws.Cells.Item(101, 2).Value = "237682117915"
$ws.Cells.Item(101, 3).Value = "MEKUEKO FOUDJO BERLINE DIDIANE ETS MOBILE FINANCIAL SERVICES MFS"
$ws.Cells.Item(101, 4).Value = 272599

$ws.Range("B86:B101").ClearFormats()
